$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("B1").Value = "X"
$ws.Range("C1").Value = "Y"

# Make the new header cell match the style (bold/border/alignment) of its neighbours
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Desvio Padrão"

# --- New data values for columns B (X), C (Y) and D (Desvio Padrao) ---
$data = @(
    @(1, 4, 0.14),
    @(2, 4.5, 0.27),
    @(3, 8, 0.055),
    @(4, 7.5, 0.195),
    @(5, 11, 0.375),
    @(6, 10.5, 0.08500000000000001),
    @(7, 14, 0.665),
    @(8, 14.5, 0.7350000000000001),
    @(9, 15, 0.46),
    @(10, 17.5, 0.25),
    @(11, 20, 0.135),
    @(12, 22.5, 0.29),
    @(13, 20, 0.31),
    @(14, 21.5, 0.495),
    @(15, 26, 1.05),
    @(16, 28.5, 1.11),
    @(17, 29, 0.39),
    @(18, 29.5, 1.435),
    @(19, 30, 1.935),
    @(20, 31.5, 2.25)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

# --- Remove the now-unused rows 22-27 (previously held extra data points) ---
$ws.Range("A22:D27").Clear()
